$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data source (EC database) was refreshed and re-sorted: the table that used
# to be grouped by employee (all periods for FABIAN, then all periods for ESTEFANY) is now
# grouped by period (ascending 2111..2206), alternating FABIAN / ESTEFANY for each period.
# The underlying (period, mora, salario) facts per employee are unchanged - only the row
# order / grouping changed. Column B (Tipo Doc) stays "CC" for every row and is left as-is.

$rows = @(
    @{ Row = 16; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2111"; Mora = 120000; Salario = 3000000 },
    @{ Row = 17; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2111"; Mora = 82609;  Salario = 2065217 },
    @{ Row = 18; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2112"; Mora = 120000; Salario = 3000000 },
    @{ Row = 19; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2112"; Mora = 82609;  Salario = 2065217 },
    @{ Row = 20; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2201"; Mora = 120000; Salario = 3000000 },
    @{ Row = 21; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2201"; Mora = 82609;  Salario = 2065217 },
    @{ Row = 22; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2202"; Mora = 120000; Salario = 3000000 },
    @{ Row = 23; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2202"; Mora = 82609;  Salario = 2065217 },
    @{ Row = 24; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2203"; Mora = 120000; Salario = 3000000 },
    @{ Row = 25; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2203"; Mora = 82609;  Salario = 2065217 },
    @{ Row = 26; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2204"; Mora = 120000; Salario = 3000000 },
    @{ Row = 27; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2204"; Mora = 82609;  Salario = 2065217 },
    @{ Row = 28; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2205"; Mora = 120000; Salario = 3000000 },
    @{ Row = 29; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2205"; Mora = 82609;  Salario = 2065217 },
    @{ Row = 30; Doc = "73127525";    Nombre = "FABIAN EUGENIO PINEDA LOPEZ";      Periodo = "2206"; Mora = 88000;  Salario = 3000000 },
    @{ Row = 31; Doc = "1128061946";  Nombre = "ESTEFANY DEL CARMEN GUZMAN AVILA"; Periodo = "2206"; Mora = 60580;  Salario = 2065217 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value2 = $r.Doc
    $ws.Range("D$n").Value2 = $r.Nombre
    $ws.Range("E$n").Value2 = $r.Periodo
    $ws.Range("F$n").Value2 = $r.Mora
    $ws.Range("G$n").Value2 = $r.Salario
}
